$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so they stay as text (matching original inlineStr type)
$textCells = @("D5", "D10", "D16", "D18", "D24", "D26", "D31", "D34", "D36", "D40", "D42", "D43", "D45", "D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.737.67'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '1.601.87'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").Value = '211.85'
$ws.Range("E5").Value = '  +0.12%  '
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("E8").Value = '  +0.18%  '
$ws.Range("E9").Value = '  +0.29%  '
$ws.Range("D10").Value = '19.78'
$ws.Range("E10").Value = '  +1.03%  '
$ws.Range("E11").Value = '  +0.63%  '
$ws.Range("D12").Value = '1.827.03'
$ws.Range("E12").Value = '  +0.17%  '
$ws.Range("D13").Value = '1.603.61'
$ws.Range("E13").Value = '  +0.41%  '
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("E15").Value = '  -0.19%  '
$ws.Range("D16").Value = '65.04'
$ws.Range("E16").Value = '  -0.14%  '
$ws.Range("D17").Value = '0.0₃0739'
$ws.Range("E17").Value = '  +0.48%  '
$ws.Range("D18").Value = '210.41'
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("E19").Value = '  +0.27%  '
$ws.Range("E20").Value = '  +2.23%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("E22").Value = '  -2.00%  '
$ws.Range("E23").Value = '  +0.35%  '
$ws.Range("D24").Value = '143.71'
$ws.Range("E24").Value = '  -0.98%  '
$ws.Range("E25").Value = '  +0.28%  '
$ws.Range("D26").Value = '7.10'
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("E27").Value = '  -0.74%  '
$ws.Range("E28").Value = '  +0.54%  '
$ws.Range("E29").Value = '  -0.46%  '
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("D31").Value = '3.28'
$ws.Range("E31").Value = '  +1.28%  '
$ws.Range("E32").Value = '  +1.18%  '
$ws.Range("D33").Value = '1.294.85'
$ws.Range("E33").Value = '  +1.45%  '
$ws.Range("D34").Value = '2.47'
$ws.Range("E34").Value = '  +0.77%  '
$ws.Range("E35").Value = '  +0.89%  '
$ws.Range("D36").Value = '0.602'
$ws.Range("E36").Value = '  -3.06%  '
$ws.Range("E37").Value = '  +11.61%  '
$ws.Range("E38").Value = '  -0.38%  '
$ws.Range("E39").Value = '  -0.56%  '
$ws.Range("D40").Value = '5.41'
$ws.Range("E40").Value = '  -1.82%  '
$ws.Range("E41").Value = '  -0.25%  '
$ws.Range("D42").Value = '0.786'
$ws.Range("E42").Value = '  +0.33%  '
$ws.Range("D43").Value = '63.13'
$ws.Range("D44").Value = '1.738.95'
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("D45").Value = '90.68'
$ws.Range("E45").Value = '  -0.39%  '
$ws.Range("E46").Value = '  -2.28%  '
$ws.Range("E47").Value = '  -0.08%  '
$ws.Range("D48").Value = '0.0517'
$ws.Range("E48").Value = '  +1.75%  '
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("E50").Value = '  +0.75%  '
$ws.Range("E51").Value = '  +1.04%  '
